# Actualizacion Datos Personales 4 nov
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("3ASV")
$ws2 = $wb.Worksheets.Item("3APM")

# --- Sheet 3ASV ---
# Row 8 (JESUS SAMUEL GARCIA LEON): correo changed
$ws1.Range("E8").Value = "erikaleonpalacios7@gmail.com"

# Row 4 (DIANA ITZEL BONILLA TEPEPA): tutor name typo fix ITSEL -> ITZEL
$ws1.Range("H4").Value = "GUADALUPE ITZEL TEPEPA ROSAS"

# Row 15 (IRVING MOLINA MORALES): tutor name duplicated text fix
$ws1.Range("H15").Value = "MIRIAM MOLINA MORALES"

# Row 19 (KELLY ITZEL RIVERA VARGAS): add missing tutor e-mail
$ws1.Range("I19").Value = "Kelly_rivera_vargas@gmail.com"

# --- Sheet 3APM ---
# Row 6 (CESAR CUEVAS CUATRA): fill in previously empty data
$ws2.Range("E6").Value = "cesarcuevasc3@gmail.com"
$ws2.Range("F6").Value = "2722848082"
$ws2.Range("G6").Value = "2722848082"
$ws2.Range("H6").Value = "GERARDO CUEVAS MACUIXTLE"
$ws2.Range("I6").Value = "cesarcuevasc3@gmail.com"
$ws2.Range("J6").Value = "2722848082"
